$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 78
$ws.Range("I2").Value = 78
$ws.Range("K2").Value = 78
$ws.Range("M2").Value = 35
# Row 6
$ws.Range("H6").Value = 7544.5884
$ws.Range("I6").Value = 1282.7142
$ws.Range("K6").Value = 3848.1426
$ws.Range("M6").Value = -3736.1426
# Row 15
$ws.Range("H15").Value = 5100.7466
$ws.Range("I15").Value = 5100.7466
$ws.Range("K15").Value = 15302.2398
$ws.Range("M15").Value = -15133.2398
# Row 45
$ws.Range("H45").Value = 6566.75
$ws.Range("J45").Value = 10000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -30384
# Row 112
$ws.Range("H112").Value = 1713.6923
$ws.Range("J112").Value = 1776.24
$ws.Range("L112").Value = 5328.72
$ws.Range("N112").Value = -7544.72
# Row 129
$ws.Range("H129").Value = 21305.734
$ws.Range("I129").Value = 773.43475
$ws.Range("J129").Value = 39468.92
$ws.Range("K129").Value = 2320.30425
$ws.Range("L129").Value = 118406.76
$ws.Range("M129").Value = 2679.69575
$ws.Range("N129").Value = -128406.76
# Row 132
$ws.Range("H132").Value = 28353.93
$ws.Range("I132").Value = 16736.656
$ws.Range("J132").Value = 222943.25
$ws.Range("K132").Value = 50209.96799999999
$ws.Range("L132").Value = 668829.75
$ws.Range("M132").Value = -47679.96799999999
$ws.Range("N132").Value = -673889.75
# Row 135
$ws.Range("H135").Value = 33334660
$ws.Range("I135").Value = 1376.3077
$ws.Range("J135").Value = 250001000
$ws.Range("K135").Value = 12386.7693
$ws.Range("L135").Value = 2250009000
$ws.Range("M135").Value = -9851.7693
$ws.Range("N135").Value = -2250014070
# Row 137
$ws.Range("H137").Value = 4490.6665
$ws.Range("I137").Value = 3060.6
$ws.Range("J137").Value = 9117.352999999999
$ws.Range("K137").Value = 9181.799999999999
$ws.Range("L137").Value = 27352.059
$ws.Range("M137").Value = -6631.799999999999
$ws.Range("N137").Value = -32452.059
# Row 138
$ws.Range("H138").Value = 1754.8103
$ws.Range("I138").Value = 1403.4138
$ws.Range("J138").Value = 2106.2068
$ws.Range("K138").Value = 4210.2414
$ws.Range("L138").Value = 6318.6204
$ws.Range("M138").Value = 929.7586000000001
$ws.Range("N138").Value = -16598.6204
# Row 141
$ws.Range("H141").Value = 2105.0908
$ws.Range("I141").Value = 977.2083
$ws.Range("K141").Value = 2931.6249
$ws.Range("M141").Value = 2248.3751

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1483.4364
$ws.Range("I61").Value = 988
$ws.Range("J61").Value = 2285.5715
$ws.Range("K61").Value = 988
$ws.Range("L61").Value = 2285.5715
$ws.Range("M61").Value = -776
$ws.Range("N61").Value = -2709.5715
# Row 132
$ws.Range("H132").Value = 23812570
$ws.Range("I132").Value = 41669890
$ws.Range("J132").Value = 2810.4443
$ws.Range("K132").Value = 125009670
$ws.Range("L132").Value = 8431.332900000001
$ws.Range("M132").Value = -125007140
$ws.Range("N132").Value = -13491.3329
# Row 136
$ws.Range("H136").Value = 1483.4364
$ws.Range("I136").Value = 988
$ws.Range("J136").Value = 2285.5715
$ws.Range("K136").Value = 2964
$ws.Range("L136").Value = 6856.7145
$ws.Range("M136").Value = -414
$ws.Range("N136").Value = -11956.7145

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# Row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
# Row 58
$ws.Range("H58").Value = 22279.5
$ws.Range("J58").Value = 22279.5
$ws.Range("L58").Value = 22279.5
$ws.Range("N58").Value = -22867.5
# Row 59
$ws.Range("H59").Value = 41266.668
$ws.Range("J59").Value = 41266.668
$ws.Range("L59").Value = 41266.668
$ws.Range("N59").Value = -42960.668
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 134
$ws.Range("H134").Value = 1978.4694
$ws.Range("I134").Value = 1698.3611
$ws.Range("K134").Value = 5095.0833
$ws.Range("M134").Value = -2560.0833

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1029.51
$ws.Range("I31").Value = 691.3393
$ws.Range("J31").Value = 1459.909
$ws.Range("K31").Value = 691.3393
$ws.Range("L31").Value = 1459.909
$ws.Range("M31").Value = -396.3393
$ws.Range("N31").Value = -2049.909
# Row 34
$ws.Range("H34").Value = 1029.51
$ws.Range("I34").Value = 691.3393
$ws.Range("J34").Value = 1459.909
$ws.Range("K34").Value = 691.3393
$ws.Range("L34").Value = 1459.909
$ws.Range("M34").Value = -489.3393
$ws.Range("N34").Value = -1863.909
# Row 58
$ws.Range("H58").Value = 2444.842
$ws.Range("I58").Value = 1533.0714
$ws.Range("J58").Value = 4997.8
$ws.Range("K58").Value = 1533.0714
$ws.Range("L58").Value = 4997.8
$ws.Range("M58").Value = -1330.0714
$ws.Range("N58").Value = -5403.8
# Row 107
$ws.Range("H107").Value = 612.93335
$ws.Range("I107").Value = 373.83334
$ws.Range("J107").Value = 971.5833
$ws.Range("K107").Value = 373.83334
$ws.Range("L107").Value = 971.5833
$ws.Range("M107").Value = 1546.16666
$ws.Range("N107").Value = -4811.5833
# Row 132
$ws.Range("H132").Value = 454247.6
$ws.Range("I132").Value = 1570.64
$ws.Range("J132").Value = 2340401.5
$ws.Range("K132").Value = 4711.92
$ws.Range("L132").Value = 7021204.5
$ws.Range("M132").Value = -2181.92
$ws.Range("N132").Value = -7026264.5
# Row 134
$ws.Range("H134").Value = 351106.3
$ws.Range("I134").Value = 932.1111
$ws.Range("J134").Value = 3502674.2
$ws.Range("K134").Value = 2796.3333
$ws.Range("L134").Value = 10508022.6
$ws.Range("M134").Value = -261.3332999999998
$ws.Range("N134").Value = -10513092.6
# Row 136
$ws.Range("H136").Value = 2444.842
$ws.Range("I136").Value = 1533.0714
$ws.Range("J136").Value = 4997.8
$ws.Range("K136").Value = 4599.2142
$ws.Range("L136").Value = 14993.4
$ws.Range("M136").Value = -2049.2142
$ws.Range("N136").Value = -20093.4
# Row 138
$ws.Range("H138").Value = 33184.367
$ws.Range("J138").Value = 33184.367
$ws.Range("L138").Value = 33184.367
$ws.Range("N138").Value = -43464.367
# Row 139
$ws.Range("H139").Value = 73666
$ws.Range("J139").Value = 73666
$ws.Range("L139").Value = 73666
$ws.Range("N139").Value = -83946

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 7684368
$ws.Range("I68").Value = 4445212
$ws.Range("J68").Value = 11365227
$ws.Range("K68").Value = 13335636
$ws.Range("L68").Value = 34095681
$ws.Range("M68").Value = -13334825
$ws.Range("N68").Value = -34097303
# Row 71
$ws.Range("H71").Value = 7684368
$ws.Range("I71").Value = 4445212
$ws.Range("J71").Value = 11365227
$ws.Range("K71").Value = 40006908
$ws.Range("L71").Value = 102287043
$ws.Range("M71").Value = -40002852
$ws.Range("N71").Value = -102295155

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1245.375
$ws.Range("I122").Value = 1009
$ws.Range("J122").Value = 2269.6667
$ws.Range("K122").Value = 3027
$ws.Range("L122").Value = 6809.000100000001
$ws.Range("M122").Value = -577
$ws.Range("N122").Value = -11709.0001
# Row 127
$ws.Range("H127").Value = 41205.453
$ws.Range("J127").Value = 41205.453
$ws.Range("L127").Value = 41205.453
$ws.Range("N127").Value = -51125.453
# Row 132
$ws.Range("H132").Value = 1808.1372
$ws.Range("I132").Value = 1139.4375
$ws.Range("K132").Value = 3418.3125
$ws.Range("M132").Value = -888.3125

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2911.4348
$ws.Range("I16").Value = 2798.2632
$ws.Range("J16").Value = 3449
$ws.Range("K16").Value = 2798.2632
$ws.Range("L16").Value = 3449
$ws.Range("M16").Value = -2628.2632
$ws.Range("N16").Value = -3789
# Row 55
$ws.Range("H55").Value = 498.72726
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 497.84616
$ws.Range("K55").Value = 500
$ws.Range("L55").Value = 497.84616
$ws.Range("M55").Value = -327
$ws.Range("N55").Value = -843.8461600000001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
# Row 132
$ws.Range("H132").Value = 1709.7593
$ws.Range("I132").Value = 1346.7428
$ws.Range("J132").Value = 2378.4736
$ws.Range("K132").Value = 4040.2284
$ws.Range("L132").Value = 7135.4208
$ws.Range("M132").Value = -1510.2284
$ws.Range("N132").Value = -12195.4208
